# Update the Estado de Cuenta data table: remove the previous periods
# and add the new periods for each worker (EDGAR LUIS ALMAGRO MENDOZA and
# NELSON LUIS PEREA ANDRADE), shifting the period/value window forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047472256"
$ws.Range("D16").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E16").Value = "2312"
$ws.Range("F16").Value = 25333
$ws.Range("C17").Value = "1047472256"
$ws.Range("D17").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E17").Value = "2311"
$ws.Range("F17").Value = 40000
$ws.Range("C18").Value = "1047472256"
$ws.Range("D18").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E18").Value = "2310"
$ws.Range("F18").Value = 40000
$ws.Range("C19").Value = "1047472256"
$ws.Range("D19").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E19").Value = "2309"
$ws.Range("F19").Value = 40000
$ws.Range("C20").Value = "1047472256"
$ws.Range("D20").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E20").Value = "2308"
$ws.Range("F20").Value = 40000
$ws.Range("C21").Value = "1047472256"
$ws.Range("D21").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E21").Value = "2307"
$ws.Range("F21").Value = 40000
$ws.Range("C22").Value = "1047472256"
$ws.Range("D22").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E22").Value = "2306"
$ws.Range("F22").Value = 40000
$ws.Range("C23").Value = "1047472256"
$ws.Range("D23").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E23").Value = "2305"
$ws.Range("F23").Value = 40000
$ws.Range("C24").Value = "1047472256"
$ws.Range("D24").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E24").Value = "2304"
$ws.Range("F24").Value = 40000
$ws.Range("C25").Value = "1047472256"
$ws.Range("D25").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E25").Value = "2303"
$ws.Range("F25").Value = 40000
$ws.Range("C26").Value = "1047472256"
$ws.Range("D26").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E26").Value = "2302"
$ws.Range("F26").Value = 40000
$ws.Range("C27").Value = "1047472256"
$ws.Range("D27").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E27").Value = "2301"
$ws.Range("F27").Value = 40000
$ws.Range("C28").Value = "1047472256"
$ws.Range("D28").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E28").Value = "2212"
$ws.Range("F28").Value = 40000
$ws.Range("C29").Value = "1047472256"
$ws.Range("D29").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E29").Value = "2211"
$ws.Range("F29").Value = 40000
$ws.Range("C30").Value = "1047472256"
$ws.Range("D30").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E30").Value = "2210"
$ws.Range("F30").Value = 40000
$ws.Range("C31").Value = "1047472256"
$ws.Range("D31").Value = "EDGAR LUIS ALMAGRO MENDOZA"
$ws.Range("E31").Value = "2209"
$ws.Range("F31").Value = 40000
$ws.Range("C32").Value = "12917341"
$ws.Range("D32").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E32").Value = "2312"
$ws.Range("F32").Value = 25333
$ws.Range("C33").Value = "12917341"
$ws.Range("D33").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E33").Value = "2311"
$ws.Range("F33").Value = 40000
$ws.Range("C34").Value = "12917341"
$ws.Range("D34").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E34").Value = "2310"
$ws.Range("F34").Value = 40000
$ws.Range("C35").Value = "12917341"
$ws.Range("D35").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E35").Value = "2309"
$ws.Range("F35").Value = 40000
$ws.Range("C36").Value = "12917341"
$ws.Range("D36").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E36").Value = "2308"
$ws.Range("F36").Value = 40000
$ws.Range("C37").Value = "12917341"
$ws.Range("D37").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E37").Value = "2307"
$ws.Range("F37").Value = 40000
$ws.Range("C38").Value = "12917341"
$ws.Range("D38").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E38").Value = "2306"
$ws.Range("F38").Value = 40000
$ws.Range("C39").Value = "12917341"
$ws.Range("D39").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E39").Value = "2305"
$ws.Range("F39").Value = 40000
$ws.Range("C40").Value = "12917341"
$ws.Range("D40").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E40").Value = "2304"
$ws.Range("F40").Value = 40000
$ws.Range("C41").Value = "12917341"
$ws.Range("D41").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E41").Value = "2303"
$ws.Range("F41").Value = 40000
$ws.Range("C42").Value = "12917341"
$ws.Range("D42").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E42").Value = "2302"
$ws.Range("F42").Value = 40000
$ws.Range("C43").Value = "12917341"
$ws.Range("D43").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E43").Value = "2301"
$ws.Range("F43").Value = 40000
$ws.Range("C44").Value = "12917341"
$ws.Range("D44").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E44").Value = "2212"
$ws.Range("F44").Value = 40000
$ws.Range("C45").Value = "12917341"
$ws.Range("D45").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E45").Value = "2211"
$ws.Range("F45").Value = 40000
$ws.Range("C46").Value = "12917341"
$ws.Range("D46").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E46").Value = "2210"
$ws.Range("F46").Value = 40000
$ws.Range("C47").Value = "12917341"
$ws.Range("D47").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E47").Value = "2209"
$ws.Range("F47").Value = 40000
$ws.Range("C48").Value = "12917341"
$ws.Range("D48").Value = "NELSON LUIS PEREA ANDRADE"
$ws.Range("E48").Value = "2208"
$ws.Range("F48").Value = 12000
